# Update Name of Algo
# Apply updated KNN-imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = -12.338
    "B3"  = 5.695
    "D3"  = -7.256
    "D12" = -7.192
    "B14" = 5.857
    "B21" = 9.643000000000001
    "B23" = 8.004999999999999
    "D24" = -7.524000000000001
    "B25" = 5.531000000000001
    "C25" = -13.499
    "D25" = -8.307
    "B26" = 6.263
    "C27" = -13.055
    "B29" = 5.372999999999999
    "C31" = -13.016
    "C39" = -12.847
    "C48" = -11.25
    "D50" = -7.92
    "C51" = -11.701
    "C52" = -11.221
    "B53" = 6.383
    "D53" = -6.994
    "C55" = -13.748
    "C56" = -12.753
    "B57" = 5.178
    "C57" = -13.716
    "D57" = -8.404999999999999
    "B59" = 4.714
    "D61" = -7.822
    "D63" = -7.522
    "B69" = 5.319999999999999
    "D70" = -7.213000000000001
    "C73" = -12.654
    "B79" = 5.771
    "B83" = 5.702000000000001
    "D86" = -8.205
    "C89" = -11.12
    "C90" = -12.91
    "B91" = 5.545
    "C92" = -11.223
    "B93" = 5.589
    "D98" = -8.486000000000001
    "D100" = -8.409000000000001
    "D102" = -7.958000000000001
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
